$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D header: "statrate18" (mirrors A1/B1/C1 = "Country code"/"statrate16"/"statrate17")
$hdr = $ws.Cells.Item(1, 4)
$hdr.Value2 = "statrate18"

# Column D (statrate18) duplicates column C (statrate17) for every data row (2-187),
# matching style (numeric, style index used by B/C columns) and value, including
# blank cells where statrate17 itself has no value.
for ($r = 2; $r -le 187; $r++) {
    $src = $ws.Cells.Item($r, 3)
    $dst = $ws.Cells.Item($r, 4)
    $dst.NumberFormat = "0"
    $dst.Value2 = $src.Value2
}

# Update selection/view state to reflect the newly populated column.
$ws.Range("D1:D187").Select()
